$d = $word.ActiveDocument

# --- 1. Remove bold from the first policy bullet's heading ---
# ("ACCESS Core Information Security Policy and Procedures")
# This affects both the paragraph-mark run properties (pPr/rPr) and the
# hyperlink's own run properties.
$p18 = $d.Paragraphs.Item(18)
$p18.Range.Font.Bold = 0

$h = $d.Hyperlinks.Item(1)
$h.Range.Font.Bold = 0

# --- 2. Append " (Under Development)" to four policy bullet texts ---
$d.Content.Find.Execute(
    "ACCESS Training and Awareness Policy", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "ACCESS Training and Awareness Policy (Under Development)", 2)

$d.Content.Find.Execute(
    "ACCESS Identity and Access Management Policy", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "ACCESS Identity and Access Management Policy (Under Development)", 2)

$d.Content.Find.Execute(
    "ACCESS Information Classification Policy", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "ACCESS Information Classification Policy (Under Development)", 2)

$d.Content.Find.Execute(
    "ACCESS Disaster Recovery Policy", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "ACCESS Disaster Recovery Policy (Under Development)", 2)
